# Daily update for 20/05/2020
# Appends the 20/05/2020 (serial 43971) row of data to the four data
# tables: Table 1 - Cumulative cases, Table 2 - ICU patients,
# Table 3a - Hospital Confirmed, Table 3b - Hospital Suspected.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Table 1 - Cumulative cases -> data fills the already-present blank row 78
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1 - Cumulative cases")

# Row 78's date cell needs the same "date" formatting as the row above it,
# and the total column needs the heavier bottom border used further up the
# sheet (matches the pre-existing P58/P59/P60 "section total" style).
$ws1.Range("A77").Copy()
$ws1.Range("A78").PasteSpecial(-4122)
$ws1.Range("P58").Copy()
$ws1.Range("P78").PasteSpecial(-4122)

$ws1.Range("A78").Value = 43971
$ws1.Range("B78").Value = 987
$ws1.Range("C78").Value = 322
$ws1.Range("D78").Value = 257
$ws1.Range("E78").Value = 826
$ws1.Range("F78").Value = 894
$ws1.Range("G78").Value = 1216
$ws1.Range("H78").Value = 3780
$ws1.Range("I78").Value = 334
$ws1.Range("J78").Value = 1867
$ws1.Range("K78").Value = 2584
$ws1.Range("L78").Value = 7
$ws1.Range("M78").Value = 54
$ws1.Range("N78").Value = 1617
$ws1.Range("O78").Value = 6
$ws1.Range("P78").Value = 14751

# ---------------------------------------------------------------------------
# Table 2 - ICU patients -> brand-new row 67
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table 2 - ICU patients")

$ws2.Range("A66:P66").Copy()
$ws2.Range("A67:P67").PasteSpecial(-4122)
$ws2.Range("Q65").Copy()
$ws2.Range("Q67").PasteSpecial(-4122)

$ws2.Range("A67").Value = 43971
$ws2.Range("B67").Value = "*"
$ws2.Range("C67").Value = "*"
$ws2.Range("D67").Value = "*"
$ws2.Range("E67").Value = "*"
$ws2.Range("F67").Value = "*"
$ws2.Range("G67").Value = 6
$ws2.Range("H67").Value = 18
$ws2.Range("I67").Value = "*"
$ws2.Range("J67").Value = 8
$ws2.Range("K67").Value = 8
$ws2.Range("L67").Value = "*"
$ws2.Range("M67").Value = "*"
$ws2.Range("N67").Value = "*"
$ws2.Range("O67").Value = "*"
$ws2.Range("P67").Value = "*"
$ws2.Range("Q67").Value = 53

# ---------------------------------------------------------------------------
# Table 3a - Hospital Confirmed -> brand-new row 59
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table 3a - Hospital Confirmed")

$ws3.Range("A58:P58").Copy()
$ws3.Range("A59:P59").PasteSpecial(-4122)
$ws2.Range("Q65").Copy()
$ws3.Range("Q59").PasteSpecial(-4122)

$ws3.Range("A59").Value = 43971
$ws3.Range("B59").Value = 29
$ws3.Range("C59").Value = 21
$ws3.Range("D59").Value = "*"
$ws3.Range("E59").Value = 72
$ws3.Range("F59").Value = 15
$ws3.Range("G59").Value = 80
$ws3.Range("H59").Value = 412
$ws3.Range("I59").Value = 13
$ws3.Range("J59").Value = 91
$ws3.Range("K59").Value = 183
$ws3.Range("L59").Value = "*"
$ws3.Range("M59").Value = "*"
$ws3.Range("N59").Value = 21
$ws3.Range("O59").Value = "*"
$ws3.Range("P59").Value = "*"
$ws3.Range("Q59").Value = 943

# ---------------------------------------------------------------------------
# Table 3b- Hospital Suspected -> brand-new row 59
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Table 3b- Hospital Suspected")

$ws4.Range("A58:P58").Copy()
$ws4.Range("A59:P59").PasteSpecial(-4122)
$ws2.Range("Q65").Copy()
$ws4.Range("Q59").PasteSpecial(-4122)

$ws4.Range("A59").Value = 43971
$ws4.Range("B59").Value = 69
$ws4.Range("C59").Value = 20
$ws4.Range("D59").Value = 22
$ws4.Range("E59").Value = 31
$ws4.Range("F59").Value = 47
$ws4.Range("G59").Value = 33
$ws4.Range("H59").Value = "N/A"
$ws4.Range("I59").Value = 33
$ws4.Range("J59").Value = 100
$ws4.Range("K59").Value = 124
$ws4.Range("L59").Value = "*"
$ws4.Range("M59").Value = "*"
$ws4.Range("N59").Value = 16
$ws4.Range("O59").Value = "*"
$ws4.Range("P59").Value = "*"
$ws4.Range("Q59").Value = 500
